$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(51,1).Value = "F049"
$ws.Cells.Item(51,2).Value = "# 50. Impoundment of unsanitary equipment or food"
$ws.Cells.Item(51,3).Value = 5928
$ws.Cells.Item(52,1).Value = "F050"
$ws.Cells.Item(52,2).Value = "# 51. Permit Suspension"
$ws.Cells.Item(52,3).Value = 2955
$ws.Cells.Item(53,1).Value = "F051"
$ws.Cells.Item(53,2).Value = "# 49. Samples Collected"
$ws.Cells.Item(53,3).Value = 40
$ws.Cells.Item(54,1).Value = "F052"
$ws.Cells.Item(54,2).Value = "# 01b. Food safety certification"
$ws.Cells.Item(54,3).Value = 18359
$ws.Cells.Item(55,1).Value = "F053"
$ws.Cells.Item(55,2).Value = "# 21a. Hot Water Available"
$ws.Cells.Item(55,3).Value = 4218
$ws.Cells.Item(56,1).Value = "F054"
$ws.Cells.Item(56,2).Value = "# 52. Multiple Major Critical Violations / Increased Risk to Public Health"
$ws.Cells.Item(56,3).Value = 1214
$ws.Cells.Item(57,1).Value = "F055"
$ws.Cells.Item(57,2).Value = "# 01a. Demonstration of knowledge"
$ws.Cells.Item(57,3).Value = 1515
$ws.Cells.Item(58,1).Value = "F056"
$ws.Cells.Item(58,2).Value = "# 10. Proper cooking time & temperatures"
$ws.Cells.Item(58,3).Value = 12
$ws.Cells.Item(59,1).Value = "F057"
$ws.Cells.Item(59,2).Value = "# 18. Compliance with variance, specialized process, & HACCP Plan"
$ws.Cells.Item(59,3).Value = 43
$ws.Cells.Item(60,1).Value = "F058"
$ws.Cells.Item(60,2).Value = "# 19. Consumer advisory provided for raw or undercooked foods"
$ws.Cells.Item(60,3).Value = 185
$ws.Cells.Item(61,1).Value = "MF07"
$ws.Cells.Item(61,2).Value = "# 07. Adequate handwashing facilities supplied & accessible"
$ws.Cells.Item(61,3).Value = 1
$ws.Cells.Item(62,1).Value = "MF42"
$ws.Cells.Item(62,2).Value = "# 42. Garbage and refuse properly disposed; facilities maintained"
$ws.Cells.Item(62,3).Value = 1
$ws.Cells.Item(63,1).Value = "SF15"
$ws.Cells.Item(63,2).Value = "No Health Code Violations Observed At The Time Of Inspection"
$ws.Cells.Item(63,3).Value = 1
$ws.Cells.Item(64,1).Value = "SS33"
$ws.Cells.Item(64,2).Value = "Garbage / rubbish receptacles not maintained clean and sanitary"
$ws.Cells.Item(64,3).Value = 1
$ws.Cells.Item(65,1).Value = "W001"
$ws.Cells.Item(65,2).Value = "Proper hot and cold holding temperatures"
$ws.Cells.Item(65,3).Value = 9
$ws.Cells.Item(66,1).Value = "W002"
$ws.Cells.Item(66,2).Value = "Food in good condition, safe and unadultered"
$ws.Cells.Item(66,3).Value = 1
$ws.Cells.Item(67,1).Value = "W003"
$ws.Cells.Item(67,2).Value = "Food storage separated and protected"
$ws.Cells.Item(67,3).Value = 6
$ws.Cells.Item(68,1).Value = "W004"
$ws.Cells.Item(68,2).Value = "Food storage space"
$ws.Cells.Item(68,3).Value = 3
$ws.Cells.Item(69,1).Value = "W005"
$ws.Cells.Item(69,2).Value = "Food elevated"
$ws.Cells.Item(69,3).Value = 11
$ws.Cells.Item(70,1).Value = "W006"
$ws.Cells.Item(70,2).Value = "Food packaging protected"
$ws.Cells.Item(70,3).Value = 2
$ws.Cells.Item(71,1).Value = "W008"
$ws.Cells.Item(71,2).Value = "Rodent"
$ws.Cells.Item(71,3).Value = 6
$ws.Cells.Item(72,1).Value = "W009"
$ws.Cells.Item(72,2).Value = "Cockroaches"
$ws.Cells.Item(72,3).Value = 42
$ws.Cells.Item(73,1).Value = "W011"
$ws.Cells.Item(73,2).Value = "Storage of materials 18 inches above the floor."
$ws.Cells.Item(73,3).Value = 2
$ws.Cells.Item(74,1).Value = "W012"
$ws.Cells.Item(74,2).Value = "Fly Breeding Material"
$ws.Cells.Item(74,3).Value = 3
$ws.Cells.Item(75,1).Value = "W014"
$ws.Cells.Item(75,2).Value = "Fly Breeding"
$ws.Cells.Item(75,3).Value = 5
$ws.Cells.Item(76,1).Value = "W016"
$ws.Cells.Item(76,2).Value = "Building rodent proof"
$ws.Cells.Item(76,3).Value = 3
$ws.Cells.Item(77,1).Value = "W017"
$ws.Cells.Item(77,2).Value = "Hot and cold water available"
$ws.Cells.Item(77,3).Value = 18
$ws.Cells.Item(78,1).Value = "W018"
$ws.Cells.Item(78,2).Value = "Waste water or sewage properly disposed or not discharged on the ground."
$ws.Cells.Item(78,3).Value = 3
$ws.Cells.Item(79,1).Value = "W019"
$ws.Cells.Item(79,2).Value = "Plumbing approved and maintained in good repair."
$ws.Cells.Item(79,3).Value = 15
$ws.Cells.Item(80,1).Value = "W020"
$ws.Cells.Item(80,2).Value = "Wall(s) maintained clean"
$ws.Cells.Item(80,3).Value = 9
$ws.Cells.Item(81,1).Value = "W021"
$ws.Cells.Item(81,2).Value = "Wall(s) maintained in good repair"
$ws.Cells.Item(81,3).Value = 6
$ws.Cells.Item(82,1).Value = "W022"
$ws.Cells.Item(82,2).Value = "Wall(s) constructed of approved material"
$ws.Cells.Item(82,3).Value = 1
$ws.Cells.Item(83,1).Value = "W023"
$ws.Cells.Item(83,2).Value = "Floor maintained clean"
$ws.Cells.Item(83,3).Value = 23
$ws.Cells.Item(84,1).Value = "W024"
$ws.Cells.Item(84,2).Value = "Floor maintained in good repair"
$ws.Cells.Item(84,3).Value = 6
$ws.Cells.Item(85,1).Value = "W025"
$ws.Cells.Item(85,2).Value = "Ceiling maintained clean"
$ws.Cells.Item(85,3).Value = 2
$ws.Cells.Item(86,1).Value = "W026"
$ws.Cells.Item(86,2).Value = "Ceiling maintained in good repair"
$ws.Cells.Item(86,3).Value = 3
$ws.Cells.Item(87,1).Value = "W027"
$ws.Cells.Item(87,2).Value = "Ceiling constructed of smooth, durable, and non-absorbent material"
$ws.Cells.Item(87,3).Value = 6
$ws.Cells.Item(88,1).Value = "W028"
$ws.Cells.Item(88,2).Value = "Toilet in good repair"
$ws.Cells.Item(88,3).Value = 4
$ws.Cells.Item(89,1).Value = "W029"
$ws.Cells.Item(89,2).Value = "Toilet maintained clean / sanitary"
$ws.Cells.Item(89,3).Value = 6
$ws.Cells.Item(90,1).Value = "W030"
$ws.Cells.Item(90,2).Value = "Hand sink in good repair"
$ws.Cells.Item(90,3).Value = 1
$ws.Cells.Item(91,1).Value = "W031"
$ws.Cells.Item(91,2).Value = "Hand sink maintained clean / sanitary"
$ws.Cells.Item(91,3).Value = 4
$ws.Cells.Item(92,1).Value = "W032"
$ws.Cells.Item(92,2).Value = "Toilet room floor / walls / ceiling in good repair"
$ws.Cells.Item(92,3).Value = 3
$ws.Cells.Item(93,1).Value = "W033"
$ws.Cells.Item(93,2).Value = "Toilet room floor / walls / ceiling clean"
$ws.Cells.Item(93,3).Value = 2
$ws.Cells.Item(94,1).Value = "W034"
$ws.Cells.Item(94,2).Value = "Toilet room with toilet paper / soap / towels / trash receptacle"
$ws.Cells.Item(94,3).Value = 4
$ws.Cells.Item(95,1).Value = "W035"
$ws.Cells.Item(95,2).Value = "Toilet room well ventilated"
$ws.Cells.Item(95,3).Value = 2
$ws.Cells.Item(96,1).Value = "W036"
$ws.Cells.Item(96,2).Value = "Toilet room well lighted"
$ws.Cells.Item(96,3).Value = 2
$ws.Cells.Item(97,1).Value = "W037"
$ws.Cells.Item(97,2).Value = "Toilet available"
$ws.Cells.Item(97,3).Value = 1
$ws.Cells.Item(98,1).Value = "W038"
$ws.Cells.Item(98,2).Value = "Hand sink available"
$ws.Cells.Item(98,3).Value = 2
$ws.Cells.Item(99,1).Value = "W039"
$ws.Cells.Item(99,2).Value = "Proper storage or use of hazardous materials"
$ws.Cells.Item(99,3).Value = 1
$ws.Cells.Item(100,1).Value = "W040"
$ws.Cells.Item(100,2).Value = "Compliance with shellfish tag requirements"
$ws.Cells.Item(100,3).Value = 2
$ws.Cells.Item(101,1).Value = "W041"
$ws.Cells.Item(101,2).Value = "Premises maintained clean and sanitary"
$ws.Cells.Item(101,3).Value = 6
$ws.Cells.Item(102,1).Value = "W042"
$ws.Cells.Item(102,2).Value = "Garbage / Rubbish receptacles approved type"
$ws.Cells.Item(102,3).Value = 8
$ws.Cells.Item(103,1).Value = "W043"
$ws.Cells.Item(103,2).Value = "Garbage / Rubbish receptacles maintained in good repair"
$ws.Cells.Item(103,3).Value = 3
$ws.Cells.Item(104,1).Value = "W044"
$ws.Cells.Item(104,2).Value = "Garbage / Rubbish receptacles maintained clean and sanitary"
$ws.Cells.Item(104,3).Value = 4
$ws.Cells.Item(105,1).Value = "W045"
$ws.Cells.Item(105,2).Value = "No unapproved sleeping accomodations"
$ws.Cells.Item(105,3).Value = 1
$ws.Cells.Item(106,1).Value = "W046"
$ws.Cells.Item(106,2).Value = "Live animals"
$ws.Cells.Item(106,3).Value = 2
$ws.Cells.Item(107,1).Value = "W047"
$ws.Cells.Item(107,2).Value = "Thermometer: available, maintained in good repair"
$ws.Cells.Item(107,3).Value = 8
$ws.Cells.Item(108,1).Value = "W048"
$ws.Cells.Item(108,2).Value = "Permits Available"
$ws.Cells.Item(108,3).Value = 122
$ws.Cells.Item(109,1).Value = "W049"
$ws.Cells.Item(109,2).Value = "Food from an approved source"
$ws.Cells.Item(109,3).Value = 5
$ws.Cells.Item(110,1).Value = "W050"
$ws.Cells.Item(110,2).Value = "Food properly labeled"
$ws.Cells.Item(110,3).Value = 15
$ws.Cells.Item(111,1).Value = "W051"
$ws.Cells.Item(111,2).Value = "Walls, Floors, Ceilings: approved, maintained clean and in good repair"
$ws.Cells.Item(111,3).Value = 32
$ws.Cells.Item(112,1).Value = "W052"
$ws.Cells.Item(112,2).Value = "Equipment, Shelving, Cabinets: approved, maintained clean and in good repair"
$ws.Cells.Item(112,3).Value = 31
$ws.Cells.Item(113,1).Value = "W053"
$ws.Cells.Item(113,2).Value = "Permit Suspension"
$ws.Cells.Item(113,3).Value = 35
$ws.Cells.Item(114,1).Value = "WP13"
$ws.Cells.Item(114,2).Value = "# 13. Disease Transmission - Carrier / Lesion / Rash"
$ws.Cells.Item(114,3).Value = 1
$ws.Cells.Item(115,1).Value = "WP15"
$ws.Cells.Item(115,2).Value = "# 15. Tobacco / Eating / Drinking / Habits / Behaviors"
$ws.Cells.Item(115,3).Value = 1
$ws.Cells.Item(116,1).Value = "WP16"
$ws.Cells.Item(116,2).Value = "# 16. Hair Restraints / Outer Garments / Nails / Rings"
$ws.Cells.Item(116,3).Value = 1
$ws.Cells.Item(117,1).Value = "WP18"
$ws.Cells.Item(117,2).Value = "# 18. Personal Hygiene"
$ws.Cells.Item(117,3).Value = 1

# Trailing blank row marker (row 118), matching source workbook layout
$ws.Rows.Item(118).Hidden = $true
$ws.Rows.Item(118).Hidden = $false
